$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}


# Row 2 - BNB
Set-TextValue "D2" "305.57"
Set-TextValue "E2" "-4.00%"

# Row 3 - OKB
Set-TextValue "D3" "37.00"
Set-TextValue "E3" "-7.15%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.086"
Set-TextValue "E4" "-1.03%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07714"
Set-TextValue "E5" "-6.11%"

# Row 6 - GateToken
Set-TextValue "D6" "4.347"
Set-TextValue "E6" "0.58%"

# Row 7 - was KuCoinToken, becomes FTXToken
Set-TextValue "B7" "FTXToken"
Set-TextValue "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.883"
Set-TextValue "E7" "-8.61%"

# Row 8 - was FTXToken, becomes KuCoinToken
Set-TextValue "B8" "KuCoinToken"
Set-TextValue "C8" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "8.175"
Set-TextValue "E8" "-2.63%"

# Row 9 - BTSEToken
Set-TextValue "D9" "3.036"
Set-TextValue "E9" "-8.84%"

# Row 10 - MXToken
Set-TextValue "D10" "0.9233"
Set-TextValue "E10" "-2.00%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-TextValue "D11" "0.1205"
Set-TextValue "E11" "-10.75%"

# Row 12 - WazirX
Set-TextValue "D12" "0.1864"
Set-TextValue "E12" "-6.54%"

# Row 13 - MandalaExchangeToken
Set-TextValue "D13" "0.08750"
Set-TextValue "E13" "-4.20%"

# Row 14 - BitrueCoin
Set-TextValue "D14" "0.03385"
Set-TextValue "E14" "-4.56%"

# Row 15
Set-TextValue "D15" "0.09698"
Set-TextValue "E15" "-1.11%"

# Row 16
Set-TextValue "D16" "0.001374"
Set-TextValue "E16" "-2.12%"

# Row 17
Set-TextValue "D17" "0.005937"
Set-TextValue "E17" "-5.59%"

# Row 18
Set-TextValue "D18" "3.602"
Set-TextValue "E18" "-2.36%"

# Row 19 - only E changes
Set-TextValue "E19" "-2.39%"

# Row 20
Set-TextValue "D20" "0.1276"
Set-TextValue "E20" "-3.51%"

# Row 21
Set-TextValue "D21" "5.014"
Set-TextValue "E21" "0.92%"

# Row 22
Set-TextValue "D22" "0.2596"
Set-TextValue "E22" "5.79%"

# Row 23
Set-TextValue "D23" "0.02106"
Set-TextValue "E23" "5,159.01%"

# Row 24
Set-TextValue "D24" "0.04325"
Set-TextValue "E24" "-0.89%"

# Row 25
Set-TextValue "D25" "0.001212"
Set-TextValue "E25" "-2.13%"

# Row 26
Set-TextValue "D26" "0.004220"
Set-TextValue "E26" "-11.94%"

# Row 27
Set-TextValue "D27" "0.0001351"
Set-TextValue "E27" "3.81%"

# Row 39
Set-TextValue "D39" "0.02187"
Set-TextValue "E39" "-6.57%"

# Row 40
Set-TextValue "D40" "0.04897"
Set-TextValue "E40" "-5.72%"

# Row 41
Set-TextValue "D41" "0.007554"
Set-TextValue "E41" "-2.67%"

# Row 42
Set-TextValue "D42" "0.009948"
Set-TextValue "E42" "0.77%"

# Row 43
Set-TextValue "D43" "0.1338"
Set-TextValue "E43" "-4.95%"

# Row 44
Set-TextValue "D44" "0.001995"
Set-TextValue "E44" "-4.25%"

# Row 45
Set-TextValue "D45" "0.009106"
Set-TextValue "E45" "-1.23%"

# Row 46
Set-TextValue "D46" "0.00006544"
Set-TextValue "E46" "-1.33%"

# Row 47 - only E changes
Set-TextValue "E47" "0.02%"

# Row 48 - only E changes
Set-TextValue "E48" "1.81%"

# Row 49
Set-TextValue "D49" "0.001302"
Set-TextValue "E49" "-23.05%"

# Row 50 - only E changes
Set-TextValue "E50" "0.02%"

# Row 51 - only E changes
Set-TextValue "E51" "0.02%"
